$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift the month headers forward by one month
$ws.Range("C1").Value2 = "marzo"
$ws.Range("D1").Value2 = "abril"
$ws.Range("E1").Value2 = "mayo"
$ws.Range("F1").Value2 = "junio"

# Update column widths (expressed in COM ColumnWidth units, which are
# offset from the raw OOXML "width" attribute by ~0.83 on this sheet)
$ws.Columns.Item(3).ColumnWidth = 10.17
$ws.Columns.Item(4).ColumnWidth = 10.17
$ws.Columns.Item(5).ColumnWidth = 9.17
$ws.Columns.Item(6).ColumnWidth = 10.17
